# Update the final version of abstract
#
# The sentence that ends the "Identifying parcels and networks" paragraph is
# split into two runs, with a (re-planted) "_GoBack" bookmark marking the
# boundary - this is the editing-cursor bookmark Word leaves behind at the
# last place text was changed, which is why it now sits at the point where
# the final clause was touched instead of at its old location (inside a
# comment bubble).

$d = $word.ActiveDocument

# Locate the split point: right after "unique parcels, " and before
# "resulting in an individualized parcellation scheme for each participant."
$splitRange = $d.Content
$found = $splitRange.Find.Execute("unique parcels, ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the point right after "unique parcels, ".
    $splitRange.Collapse(0)

    # Re-plant the "_GoBack" bookmark at the new edit location.
    $d.Bookmarks.Add("_GoBack", $splitRange)
}

# Word also lazily adds the built-in "FollowedHyperlink" character style to
# the style sheet the first time it is needed (e.g. after the author
# followed/visited the reference link while proofreading the final text).
$followed = $d.Styles.Add("FollowedHyperlink", 2)
$followed.BaseStyle = "DefaultParagraphFont"
$followed.Priority = 99
$followed.UnhideWhenUsed = $true
$followed.Font.Color = 7491477
$followed.Font.Underline = 1
